$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Delete()

# Row 15 deletion drops the Excel engine's tracking of the trivial
# single-cell "Checklist" merge area (A42/B42 before the delete, which
# should become A41/B41 after everything shifts up by one row). Restore
# those self-merges explicitly.
$ws.Range("A41").Merge() | Out-Null
$ws.Range("B41").Merge() | Out-Null
